$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
2896.016791289029,
3927.658649050765,
4534.029129713806,
5421.269083827,
6030.559498840906,
6525.100025618548,
6952.981220056842,
7484.473598119846,
7843.293491568781,
8290.187337687727,
8685.730032810463,
9001.70057840767,
9226.559212791988,
9552.597441704291,
9856.564911282545,
10041.54945581178,
10237.32175682418,
10480.38340354335,
10730.56883722468,
10887.0411563089,
11162.78405101858,
11285.41736594259,
11585.38660893686,
11687.4504678139,
11878.36887432134,
12114.49714478112,
12232.45668534909,
12336.5431969694,
12528.32671445921,
12638.49012198292,
12876.23976375084,
12994.18756820729,
13137.60954399211,
13387.23539579318,
13310.86209020602,
13534.82661666705,
13667.90292984921,
13802.71209478111,
13959.33622138699,
14032.85608012556,
14038.11226015806,
14098.93460621198,
14173.67659824644,
14278.4944349387,
14502.99238362768,
14794.52154245858,
14915.06365697691,
14947.84655382285,
15029.99038520521,
15134.43046697146,
15261.4251776523,
15316.65190593453,
15378.2451341638,
15502.70906282871,
15454.36020977322,
15570.32866329254,
15596.31859954413,
15606.31128178961,
15726.93835942123,
15739.24396771388,
15741.46777860923
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
